# Append rows 206:217 (A=204..215, B=decimal values) to Sheet1, matching
# the formatting already used by the existing data rows (e.g. A205/B205).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows to append: Excel row number, A value (index), B value (ratio)
$newRows = @(
    @{ Row = 206; A = 204; B = 0.1891891891891891 },
    @{ Row = 207; A = 205; B = 0.4324324324324323 },
    @{ Row = 208; A = 206; B = 0.5783783783783784 },
    @{ Row = 209; A = 207; B = 0.4808558558558558 },
    @{ Row = 210; A = 208; B = 0.4374999999999999 },
    @{ Row = 211; A = 209; B = 0.609073359073359 },
    @{ Row = 212; A = 210; B = 0.4172297297297297 },
    @{ Row = 213; A = 211; B = 0.5175675675675675 },
    @{ Row = 214; A = 212; B = 0.4324324324324323 },
    @{ Row = 215; A = 213; B = 0.2837837837837837 },
    @{ Row = 216; A = 214; B = 0.4324324324324323 },
    @{ Row = 217; A = 215; B = 0.4324324324324323 }
)

# Use the formatting of the last existing data row (A205) as the template
# for the new "A" (index) column cells, which use style index 1 (bold,
# bordered, centered) in the original workbook.
$formatSource = $ws.Range("A205")

foreach ($item in $newRows) {
    $aCell = $ws.Cells.Item($item.Row, 1)
    $bCell = $ws.Cells.Item($item.Row, 2)

    $aCell.Value = $item.A
    $bCell.Value = $item.B
}

# Copy the cell format (font/border/alignment) from A205 onto the newly
# added A206:A217 cells without touching the values already written.
$formatSource.Copy()
$ws.Range("A206:A217").PasteSpecial(-4122)
$excel.CutCopyMode = $false
